# The commit swaps the two embedded theme parts: the theme that used to
# live in ppt/theme/theme1.xml ("Integral", used by the one slide master /
# all slides) ends up with the color scheme that used to live in
# ppt/theme/theme2.xml ("Office Theme", used by the notes master), and
# vice versa. The font scheme and format scheme are byte-identical between
# the two themes, so the only real difference is the 12-slot theme color
# scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# The PowerPoint object model lets us rewrite the 12 theme colors on the
# slide master's Theme via ThemeColorScheme.Colors(i).RGB, so we push the
# "Office Theme" palette onto the deck's single Theme (theme1.xml), which
# is what every slide actually renders with.

function ConvertTo-PpRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme color scheme, in ThemeColorScheme.Colors() index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $tcs.Colors($i).RGB = ConvertTo-PpRgb $officeThemeColors[$i - 1]
}
